{"js": "// Arabic translations for the \"English\" source section of the\n// [TEMPLATE] Affiliate email - invite to seminar document.\n//\n// The document contains the same email template repeated once per\n// language (English, Portuguese, French, Thai, Vietnamese, Spanish).\n// This change only replaces strings that are still in their original\n// English within the \"English\" master section (plus a couple of\n// leftover/English strings that slipped into later sections), turning\n// them into Arabic, matching the Crowdin translation commit.\n\nconst body = context.document.body;\n\n// Helper: find the Nth (0-based) occurrence of `searchText` in the\n// document body and replace its text in place, preserving the run's\n// own formatting (bold/rtl/color/hyperlink/etc.) and xml:space.\nasync function replaceOccurrence(searchText, index, newText) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length <= index) {\n    throw new Error(\n      `Expected at least ${index + 1} match(es) of \"${searchText}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[index].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Helper: replace every occurrence of `searchText` with `newText`.\nasync function replaceAllOccurrences(searchText, newText) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1) Language label paragraph (\"English\" -> \"\u0627\u0644\u0625\u0646\u062c\u0644\u064a\u0632\u064a\u0629\").\n//    \"English\" also appears as the nav hyperlink text (1st match) and\n//    as \"Back to English\" in other language sections (matches 2-6), so\n//    target only the 2nd match (0-based index 1), the standalone\n//    section heading.\nawait replaceOccurrence(\"English\", 1, \"\u0627\u0644\u0625\u0646\u062c\u0644\u064a\u0632\u064a\u0629\");\n\n// 2) \"Brief\" label (bold run right before the \":\" in the brief table).\nawait replaceOccurrence(\"Brief\", 0, \"\u0627\u0644\u0645\u0636\u0645\u0648\u0646\");\n\n// 3) Brief description paragraph.\nawait replaceOccurrence(\n  \"An email to partners in the the target country to invite them for a one-day seminar. It will be sent via customer.io\",\n  0,\n  \"An email to partners in the the target country to invite them for a one-day seminar. \u0633\u064a\u062a\u0645 \u0625\u0631\u0633\u0627\u0644\u0647\u0627 \u0639\u0628\u0631 customer.io\"\n);\n\n// 4) \"Target audience\" label.\nawait replaceOccurrence(\"Target audience\", 0, \"\u0627\u0644\u062c\u0645\u0647\u0648\u0631 \u0627\u0644\u0645\u0633\u062a\u0647\u062f\u0641\");\n\n// 5) Email subject / heading line.\nawait replaceOccurrence(\n  \"You\\u2019re invited to our Deriv Partner Seminar\",\n  0,\n  \"\u0623\u0646\u062a \u0645\u062f\u0639\u0648 \u0625\u0644\u0649 \u0646\u062f\u0648\u0629 \u0634\u0631\u0643\u0627\u0621 Deriv\"\n);\n\n// 6) Intro paragraph.\nawait replaceOccurrence(\n  \"We\\u2019re excited to let you know that the Deriv Affiliate team will be in [CITY] in [MONTH] to meet with you, our valued partners!\",\n  0,\n  \"\u0646\u062d\u0646 \u0645\u062a\u062d\u0645\u0633\u0648\u0646 \u0644\u0625\u0639\u0644\u0627\u0645\u0643 \u0628\u0623\u0646 \u0641\u0631\u064a\u0642 Deriv Affiliate \u0633\u064a\u0643\u0648\u0646 [CITY] \u0645\u0648\u062c\u0648\u062f\u064b\u0627 [MONTH] \u0644\u0645\u0642\u0627\u0628\u0644\u062a\u0643\u060c \u0645\u0639 \u0634\u0631\u0643\u0627\u0626\u0646\u0627 \u0627\u0644\u0643\u0631\u0627\u0645!\"\n);\n\n// 7) Seminar description paragraph.\nawait replaceOccurrence(\n  \"In this one-day seminar, we\\u2019ll be providing technical and marketing support, offering the opportunity to network with other partners over a delicious lunch as well as listening to your feedback about our partnership programmes. This is your chance to get your voice heard, which will help us plan future efforts to support you better. \",\n  0,\n  \"\u0641\u064a \u0647\u0630\u0647 \u0627\u0644\u0646\u062f\u0648\u0629 \u0627\u0644\u062a\u064a \u062a\u0633\u062a\u063a\u0631\u0642 \u064a\u0648\u0645\u064b\u0627 \u0648\u0627\u062d\u062f\u064b\u0627\u060c \u0633\u0646\u0642\u062f\u0645 \u0627\u0644\u062f\u0639\u0645 \u0627\u0644\u0641\u0646\u064a \u0648\u0627\u0644\u062a\u0633\u0648\u064a\u0642\u064a\u060c \u0648\u0646\u0648\u0641\u0631 \u0627\u0644\u0641\u0631\u0635\u0629 \u0644\u0644\u062a\u0648\u0627\u0635\u0644 \u0645\u0639 \u0634\u0631\u0643\u0627\u0621 \u0622\u062e\u0631\u064a\u0646 \u062e\u0644\u0627\u0644 \u062a\u0646\u0627\u0648\u0644 \u063a\u062f\u0627\u0621 \u0637\u064a\u0628 \u0628\u0627\u0644\u0625\u0636\u0627\u0641\u0629 \u0625\u0644\u0649 \u0627\u0644\u0627\u0633\u062a\u0645\u0627\u0639 \u0625\u0644\u0649 \u0645\u0644\u0627\u062d\u0638\u0627\u062a\u0643 \u062d\u0648\u0644 \u0628\u0631\u0627\u0645\u062c \u0627\u0644\u0634\u0631\u0627\u0643\u0629 \u0644\u062f\u064a\u0646\u0627. \u0647\u0630\u0647 \u0647\u064a \u0641\u0631\u0635\u062a\u0643 \u0644\u062c\u0639\u0644 \u0635\u0648\u062a\u0643 \u0645\u0633\u0645\u0648\u0639\u064b\u0627\u060c \u0645\u0645\u0627 \u0633\u064a\u0633\u0627\u0639\u062f\u0646\u0627 \u0639\u0644\u0649 \u062a\u062e\u0637\u064a\u0637 \u0627\u0644\u062c\u0647\u0648\u062f \u0627\u0644\u0645\u0633\u062a\u0642\u0628\u0644\u064a\u0629 \u0644\u062f\u0639\u0645\u0643 \u0628\u0634\u0643\u0644 \u0623\u0641\u0636\u0644. \"\n);\n\n// 8) RSVP sentence tail (shares the paragraph with a bold \"[DATE]\" run).\nawait replaceOccurrence(\n  \". Please note that attendance is confirmed on a first come, first served basis. We look forward to seeing you there!\",\n  0,\n  \". Please note that attendance is confirmed on a first come, first served basis. \u0646\u062a\u0637\u0644\u0639 \u0625\u0644\u0649 \u0631\u0624\u064a\u062a\u0643 \u0647\u0646\u0627\u0643!\"\n);\n\n// 9) \"Send my details\" button text. Appears twice verbatim in the\n//    document (once in the English section, once left untranslated in\n//    the French section) - both are translated to Arabic.\nawait replaceAllOccurrences(\"Send my details\", \"\u0623\u0631\u0633\u0644 \u0627\u0644\u062a\u0641\u0627\u0635\u064a\u0644 \u0627\u0644\u062e\u0627\u0635\u0629 \u0628\u064a\");\n\n// 10) Contact intro sentence (English section only).\nawait replaceOccurrence(\n  \"If you have any questions, please contact us via \",\n  0,\n  \"\u0625\u0630\u0627 \u0643\u0627\u0646\u062a \u0644\u062f\u064a\u0643 \u0623\u064a \u0623\u0633\u0626\u0644\u0629\u060c \u0641\u0627\u062a\u0635\u0644 \u0628\u0646\u0627:  \"\n);\n\n// 11) \"live chat\" hyperlink text. Appears once per language section\n//     (5 total); only the English section's (1st match) changes.\nawait replaceOccurrence(\"live chat\", 0, \"\u0627\u0644\u062f\u0631\u062f\u0634\u0629 \u0627\u0644\u062d\u064a\u0629\");\n", "ps1": "# Arabic translations for the \"English\" source section of the\n# [TEMPLATE] Affiliate email - invite to seminar document.\n#\n# The document contains the same email template repeated once per\n# language (English, Portuguese, French, Thai, Vietnamese, Spanish).\n# This change only replaces strings that are still in their original\n# English within the \"English\" master section (plus a couple of\n# leftover/English strings that slipped into later sections), turning\n# them into Arabic, matching the Crowdin translation commit.\n\n$d = $word.ActiveDocument\n\n# Replace the Nth (1-based) occurrence of $searchText in the whole\n# document with $replaceText, preserving that run's own formatting\n# (bold/rtl/color/hyperlink/etc.) since only the range's .Text is set.\nfunction Replace-NthOccurrence {\n    param($doc, $searchText, $occurrenceIndex, $replaceText)\n\n    $rng = $doc.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 0\n    $rng.Find.Text = $searchText\n\n    $count = 0\n    while ($rng.Find.Execute()) {\n        $count = $count + 1\n        if ($count -eq $occurrenceIndex) {\n            $rng.Text = $replaceText\n            return $true\n        }\n        $rng.Collapse(0)\n    }\n    return $false\n}\n\n# Replace every occurrence of $searchText in the whole document with\n# $replaceText.\nfunction Replace-AllOccurrences {\n    param($doc, $searchText, $replaceText)\n\n    $rng = $doc.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 0\n    $rng.Find.Text = $searchText\n    $rng.Find.Replacement.Text = $replaceText\n    # wdReplaceAll = 2\n    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 0, $false, $replaceText, 2)\n}\n\n# 1) Language label paragraph (\"English\" -> \"\u0627\u0644\u0625\u0646\u062c\u0644\u064a\u0632\u064a\u0629\").\n#    \"English\" also appears as the nav hyperlink text (1st match) and\n#    as \"Back to English\" in other language sections (matches 3-7), so\n#    target only the 2nd match, the standalone section heading.\nReplace-NthOccurrence $d \"English\" 2 \"\u0627\u0644\u0625\u0646\u062c\u0644\u064a\u0632\u064a\u0629\" | Out-Null\n\n# 2) \"Brief\" label (bold run right before the \":\" in the brief table).\nReplace-NthOccurrence $d \"Brief\" 1 \"\u0627\u0644\u0645\u0636\u0645\u0648\u0646\" | Out-Null\n\n# 3) Brief description paragraph.\nReplace-NthOccurrence $d \"An email to partners in the the target country to invite them for a one-day seminar. It will be sent via customer.io\" 1 \"An email to partners in the the target country to invite them for a one-day seminar. \u0633\u064a\u062a\u0645 \u0625\u0631\u0633\u0627\u0644\u0647\u0627 \u0639\u0628\u0631 customer.io\" | Out-Null\n\n# 4) \"Target audience\" label.\nReplace-NthOccurrence $d \"Target audience\" 1 \"\u0627\u0644\u062c\u0645\u0647\u0648\u0631 \u0627\u0644\u0645\u0633\u062a\u0647\u062f\u0641\" | Out-Null\n\n# 5) Email subject / heading line.\nReplace-NthOccurrence $d \"You\u2019re invited to our Deriv Partner Seminar\" 1 \"\u0623\u0646\u062a \u0645\u062f\u0639\u0648 \u0625\u0644\u0649 \u0646\u062f\u0648\u0629 \u0634\u0631\u0643\u0627\u0621 Deriv\" | Out-Null\n\n# 6) Intro paragraph.\nReplace-NthOccurrence $d \"We\u2019re excited to let you know that the Deriv Affiliate team will be in [CITY] in [MONTH] to meet with you, our valued partners!\" 1 \"\u0646\u062d\u0646 \u0645\u062a\u062d\u0645\u0633\u0648\u0646 \u0644\u0625\u0639\u0644\u0627\u0645\u0643 \u0628\u0623\u0646 \u0641\u0631\u064a\u0642 Deriv Affiliate \u0633\u064a\u0643\u0648\u0646 [CITY] \u0645\u0648\u062c\u0648\u062f\u064b\u0627 [MONTH] \u0644\u0645\u0642\u0627\u0628\u0644\u062a\u0643\u060c \u0645\u0639 \u0634\u0631\u0643\u0627\u0626\u0646\u0627 \u0627\u0644\u0643\u0631\u0627\u0645!\" | Out-Null\n\n# 7) Seminar description paragraph.\nReplace-NthOccurrence $d \"In this one-day seminar, we\u2019ll be providing technical and marketing support, offering the opportunity to network with other partners over a delicious lunch as well as listening to your feedback about our partnership programmes. This is your chance to get your voice heard, which will help us plan future efforts to support you better. \" 1 \"\u0641\u064a \u0647\u0630\u0647 \u0627\u0644\u0646\u062f\u0648\u0629 \u0627\u0644\u062a\u064a \u062a\u0633\u062a\u063a\u0631\u0642 \u064a\u0648\u0645\u064b\u0627 \u0648\u0627\u062d\u062f\u064b\u0627\u060c \u0633\u0646\u0642\u062f\u0645 \u0627\u0644\u062f\u0639\u0645 \u0627\u0644\u0641\u0646\u064a \u0648\u0627\u0644\u062a\u0633\u0648\u064a\u0642\u064a\u060c \u0648\u0646\u0648\u0641\u0631 \u0627\u0644\u0641\u0631\u0635\u0629 \u0644\u0644\u062a\u0648\u0627\u0635\u0644 \u0645\u0639 \u0634\u0631\u0643\u0627\u0621 \u0622\u062e\u0631\u064a\u0646 \u062e\u0644\u0627\u0644 \u062a\u0646\u0627\u0648\u0644 \u063a\u062f\u0627\u0621 \u0637\u064a\u0628 \u0628\u0627\u0644\u0625\u0636\u0627\u0641\u0629 \u0625\u0644\u0649 \u0627\u0644\u0627\u0633\u062a\u0645\u0627\u0639 \u0625\u0644\u0649 \u0645\u0644\u0627\u062d\u0638\u0627\u062a\u0643 \u062d\u0648\u0644 \u0628\u0631\u0627\u0645\u062c \u0627\u0644\u0634\u0631\u0627\u0643\u0629 \u0644\u062f\u064a\u0646\u0627. \u0647\u0630\u0647 \u0647\u064a \u0641\u0631\u0635\u062a\u0643 \u0644\u062c\u0639\u0644 \u0635\u0648\u062a\u0643 \u0645\u0633\u0645\u0648\u0639\u064b\u0627\u060c \u0645\u0645\u0627 \u0633\u064a\u0633\u0627\u0639\u062f\u0646\u0627 \u0639\u0644\u0649 \u062a\u062e\u0637\u064a\u0637 \u0627\u0644\u062c\u0647\u0648\u062f \u0627\u0644\u0645\u0633\u062a\u0642\u0628\u0644\u064a\u0629 \u0644\u062f\u0639\u0645\u0643 \u0628\u0634\u0643\u0644 \u0623\u0641\u0636\u0644. \" | Out-Null\n\n# 8) RSVP sentence tail (shares the paragraph with a bold \"[DATE]\" run).\nReplace-NthOccurrence $d \". Please note that attendance is confirmed on a first come, first served basis. We look forward to seeing you there!\" 1 \". Please note that attendance is confirmed on a first come, first served basis. \u0646\u062a\u0637\u0644\u0639 \u0625\u0644\u0649 \u0631\u0624\u064a\u062a\u0643 \u0647\u0646\u0627\u0643!\" | Out-Null\n\n# 9) \"Send my details\" button text. Appears twice verbatim in the\n#    document (once in the English section, once left untranslated in\n#    the French section) - both are translated to Arabic.\nReplace-AllOccurrences $d \"Send my details\" \"\u0623\u0631\u0633\u0644 \u0627\u0644\u062a\u0641\u0627\u0635\u064a\u0644 \u0627\u0644\u062e\u0627\u0635\u0629 \u0628\u064a\" | Out-Null\n\n# 10) Contact intro sentence (English section only).\nReplace-NthOccurrence $d \"If you have any questions, please contact us via \" 1 \"\u0625\u0630\u0627 \u0643\u0627\u0646\u062a \u0644\u062f\u064a\u0643 \u0623\u064a \u0623\u0633\u0626\u0644\u0629\u060c \u0641\u0627\u062a\u0635\u0644 \u0628\u0646\u0627:  \" | Out-Null\n\n# 11) \"live chat\" hyperlink text. Appears once per language section\n#     (5 total); only the English section's (1st match) changes.\nReplace-NthOccurrence $d \"live chat\" 1 \"\u0627\u0644\u062f\u0631\u062f\u0634\u0629 \u0627\u0644\u062d\u064a\u0629\" | Out-Null\n"}
